# Apply the "add to openpyxl notebook" edit:
#  - replace the placeholder A1:A10 number column with a names/means/stds
#    table (A1:C7)
#  - point the chart's single series at the new "means" column and rename
#    it, and retitle the chart
#  - reposition the chart anchor to column E (idx 4), row 1 (idx 0)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- worksheet data --------------------------------------------------
$ws.Range("A1").Value = "names"
$ws.Range("B1").Value = "means"
$ws.Range("C1").Value = "stds"

$ws.Range("A2").Value = "kkp"
$ws.Range("B2").Value = 0.1043081461209843
$ws.Range("C2").Value = 0.2601834427239736

$ws.Range("A3").Value = "srl"
$ws.Range("B3").Value = 0.03331519044931429
$ws.Range("C3").Value = 0.4971198541676912

$ws.Range("A4").Value = "iok"
$ws.Range("B4").Value = 0.05019411593986726
$ws.Range("C4").Value = 2.669450422584103

$ws.Range("A5").Value = "nfm"
$ws.Range("B5").Value = 0.01011441930305049
$ws.Range("C5").Value = 0.7709364009354619

$ws.Range("A6").Value = "uwz"
$ws.Range("B6").Value = 0.1067485493249866
$ws.Range("C6").Value = 1.255990664514469

$ws.Range("A7").Value = "rjw"
$ws.Range("B7").Value = 0.05395380258988404
$ws.Range("C7").Value = 0.713453214246213

# drop the old tail of the placeholder column (rows 8-10 are gone now)
$ws.Range("A8:A10").Clear()

# --- chart -------------------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart

$chart.ChartTitle.Text = "Cytokine array results"
$chart.SeriesCollection(1).Formula = "=SERIES(""means"",,'Sheet'!`$B`$2:`$B`$8,1)"

# move the chart anchor from col C/row 5 (idx 2,4) to col E/row 1 (idx 4,0)
$chartObj.Top = $ws.Cells.Item(1, 1).Top
$chartObj.Left = $ws.Cells.Item(1, 5).Left
